$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fix the lang_code column (E): rows 11-19 are Arabic rows that were
# incorrectly tagged "eng" -> "ara"; rows 20-28 are French rows that were
# incorrectly tagged "eng" -> "fra".
# ---------------------------------------------------------------------------
$ws.Range("E11:E19").Value = "ara"
$ws.Range("E20:E28").Value = "fra"

# ---------------------------------------------------------------------------
# Normalize the two Arabic "name"/"descr" strings that had stray double
# quotes baked into the text.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"
$ws.Range("C17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"

$ws.Range("B19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"
$ws.Range("C19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"

# ---------------------------------------------------------------------------
# Column widths: B -> 41 characters, C -> 45.81640625 characters.
# (ColumnWidth is expressed in "characters"; the stored XML width includes a
# constant ~5/6 character padding offset, so we compensate for it here.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 40.166666666666664
$ws.Columns.Item(3).ColumnWidth = 44.983072916666664

# ---------------------------------------------------------------------------
# View state: scroll so row 7 is at the top, and leave C19 as the active
# selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()

# ---------------------------------------------------------------------------
# Page setup: A4, portrait.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup()
$ps.PaperSize = 9
$ps.Orientation = 1
